$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns - copy formatting from the H1 header
# cell (same bold/border/centered style used by the existing headers), then
# overwrite with the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I is a constant 1, J mirrors column H
for ($r = 2; $r -le 33; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValue
}
